$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update row 6 (ano = 2025) with refreshed metrics
$ws.Range("C6").Value = 402
$ws.Range("E6").Value = 95
$ws.Range("G6").Value = 23.6318407960199
$ws.Range("H6").Value = 76.3681592039801
